# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (want-to-go count) and "最低票价" (lowest price)
# columns (F and G) across the 展览, 演出 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3733
$ws1.Range("G4").Value = 85
$ws1.Range("F5").Value = 3733
$ws1.Range("G5").Value = 85
$ws1.Range("G6").Value = 80
$ws1.Range("F7").Value = 5270
$ws1.Range("G7").Value = 80
$ws1.Range("F8").Value = 579
$ws1.Range("F9").Value = 411
$ws1.Range("F11").Value = 1043
$ws1.Range("F22").Value = 6031
$ws1.Range("F26").Value = 6793
$ws1.Range("F30").Value = 364
$ws1.Range("F31").Value = 745
$ws1.Range("F32").Value = 4456
$ws1.Range("F36").Value = 1126
$ws1.Range("F40").Value = 914
$ws1.Range("F41").Value = 1111

# ---- Sheet "演出" (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 65

# ---- Sheet "全部类型" (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3733
$ws4.Range("G7").Value = 85
$ws4.Range("F8").Value = 3733
$ws4.Range("G8").Value = 85
$ws4.Range("G9").Value = 80
$ws4.Range("F10").Value = 5270
$ws4.Range("G10").Value = 80
$ws4.Range("F11").Value = 579
$ws4.Range("F12").Value = 411
$ws4.Range("F14").Value = 1043
$ws4.Range("F26").Value = 6031
$ws4.Range("F30").Value = 6793
$ws4.Range("F34").Value = 364
$ws4.Range("F35").Value = 745
$ws4.Range("F36").Value = 4456
$ws4.Range("F41").Value = 1126
$ws4.Range("F45").Value = 914
$ws4.Range("F46").Value = 1111
$ws4.Range("F50").Value = 65
